# Insert a new price record at row 50 ("Región del Maule", 2023-11-21),
# shifting every subsequent record (old rows 50-152) down by one row
# (new rows 51-153). This matches the commit "Fruta / hortaliza, semanal"
# which adds the latest weekly quote to the top of this block of records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 50..152 down to 51..153, leaving a blank row 50 (inherits the
# date-format style from the row that used to be there, via Excel's
# insert-shift behaviour).
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new record's data.
$ws.Cells.Item(50, 1).Value = 5
$ws.Cells.Item(50, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(50, 3).Value = "Maule"
$ws.Cells.Item(50, 4).Value = 45251
$ws.Cells.Item(50, 5).Value = 7
$ws.Cells.Item(50, 6).Value = 100112026
$ws.Cells.Item(50, 7).Value = "Haba"
$ws.Cells.Item(50, 8).Value = "Sin especificar"
$ws.Cells.Item(50, 9).Value = "Primera"
$ws.Cells.Item(50, 10).Value = 500
$ws.Cells.Item(50, 11).Value = 10000
$ws.Cells.Item(50, 12).Value = 10000
$ws.Cells.Item(50, 13).Value = 10000
$ws.Cells.Item(50, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(50, 15).Value = "Región del Maule"
$ws.Cells.Item(50, 16).Value = 400
$ws.Cells.Item(50, 17).Value = 25
$ws.Cells.Item(50, 18).Value = "Hortaliza"
